# Updated symbol list on Mon Dec 12 08:39:10 UTC 2022 with GitHub Actions
#
# This script re-applies the "Price" column refresh plus the CEJI/KickToken
# row swap observed in the diff. All D-column (Price) cells are stored as
# text in the workbook (t="inlineStr"/text cells, not numbers), so we write
# them with a leading apostrophe to force Excel to keep them as text while
# still showing/storing the plain numeric-looking string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
}

# --- Price (column D) updates -------------------------------------------------
Set-PriceText "D2"  "281.29"
Set-PriceText "D3"  "20.79"
Set-PriceText "D4"  "6.233"
Set-PriceText "D5"  "0.06140"
Set-PriceText "D6"  "3.572"
Set-PriceText "D7"  "6.561"
Set-PriceText "D8"  "1.473"
Set-PriceText "D9"  "0.8169"
Set-PriceText "D10" "0.01378"
Set-PriceText "D11" "0.1627"
Set-PriceText "D12" "0.08290"
Set-PriceText "D13" "0.03535"
Set-PriceText "D15" "0.09142"
Set-PriceText "D16" "3.710"
Set-PriceText "D17" "0.001641"
Set-PriceText "D18" "0.04646"
Set-PriceText "D19" "0.006409"
Set-PriceText "D20" "0.006165"
Set-PriceText "D21" "0.001066"
Set-PriceText "D22" "0.0001501"
Set-PriceText "D23" "3.803"
Set-PriceText "D24" "2.293"
Set-PriceText "D25" "0.3334"
Set-PriceText "D26" "0.1250"
Set-PriceText "D40" "0.04673"

# --- Row 41 / 42 swap: CEJI <-> KickToken --------------------------------------
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-PriceText "D41" "0.007174"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-PriceText "D42" "0.004002"
$ws.Range("E42").Value = "41CEJICEJI"

# --- Remaining Price (column D) updates ---------------------------------------
Set-PriceText "D43" "0.1098"
Set-PriceText "D44" "0.01132"
Set-PriceText "D45" "0.00006282"
Set-PriceText "D46" "0.00000000750"
Set-PriceText "D47" "1.001"
Set-PriceText "D48" "0.002943"
Set-PriceText "D49" "0.00001901"
Set-PriceText "D50" "0.01241"
